# Auto-generated script to update TPM-derived NATMI metrics in the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 1.458056666666667
$ws.Cells.Item(2, 8).Value = 4.374169999999999
$ws.Cells.Item(2, 9).Value = 0.2323568509805328
$ws.Cells.Item(2, 10).Value = 0.2323568509805327
$ws.Cells.Item(2, 13).Value = 29.01761566666667
$ws.Cells.Item(2, 14).Value = 87.052847
$ws.Cells.Item(2, 15).Value = 0.6436075952942075
$ws.Cells.Item(2, 16).Value = 0.6436075952942075
$ws.Cells.Item(2, 17).Value = 42.30932797355445
$ws.Cells.Item(2, 18).Value = 380.78395176199
$ws.Cells.Item(2, 19).Value = 0.1495466341097152
$ws.Cells.Item(2, 20).Value = 0.1495466341097152
$ws.Cells.Item(3, 7).Value = 1.458056666666667
$ws.Cells.Item(3, 8).Value = 4.374169999999999
$ws.Cells.Item(3, 9).Value = 0.2323568509805328
$ws.Cells.Item(3, 10).Value = 0.2323568509805327
$ws.Cells.Item(3, 15).Value = 0.05796185537580412
$ws.Cells.Item(3, 16).Value = 0.05796185537580412
$ws.Cells.Item(3, 17).Value = 3.810283108808889
$ws.Cells.Item(3, 18).Value = 34.29254797927999
$ws.Cells.Item(3, 19).Value = 0.01346783419211091
$ws.Cells.Item(3, 20).Value = 0.01346783419211091
$ws.Cells.Item(4, 7).Value = 1.458056666666667
$ws.Cells.Item(4, 8).Value = 4.374169999999999
$ws.Cells.Item(4, 9).Value = 0.2323568509805328
$ws.Cells.Item(4, 10).Value = 0.2323568509805327
$ws.Cells.Item(4, 13).Value = 1.123006
$ws.Cells.Item(4, 14).Value = 3.369018
$ws.Cells.Item(4, 15).Value = 0.02490815232594174
$ws.Cells.Item(4, 16).Value = 0.02490815232594174
$ws.Cells.Item(4, 17).Value = 1.637406385006666
$ws.Cells.Item(4, 18).Value = 14.73665746506
$ws.Cells.Item(4, 19).Value = 0.005787579838199256
$ws.Cells.Item(4, 20).Value = 0.005787579838199254
$ws.Cells.Item(5, 7).Value = 1.458056666666667
$ws.Cells.Item(5, 8).Value = 4.374169999999999
$ws.Cells.Item(5, 9).Value = 0.2323568509805328
$ws.Cells.Item(5, 10).Value = 0.2323568509805327
$ws.Cells.Item(5, 13).Value = 10.33196133333333
$ws.Cells.Item(5, 14).Value = 30.995884
$ws.Cells.Item(5, 15).Value = 0.2291617914030796
$ws.Cells.Item(5, 16).Value = 0.2291617914030796
$ws.Cells.Item(5, 17).Value = 15.06458510180889
$ws.Cells.Item(5, 18).Value = 135.58126591628
$ws.Cells.Item(5, 19).Value = 0.05324731221547729
$ws.Cells.Item(5, 20).Value = 0.05324731221547729
$ws.Cells.Item(6, 7).Value = 1.458056666666667
$ws.Cells.Item(6, 8).Value = 4.374169999999999
$ws.Cells.Item(6, 9).Value = 0.2323568509805328
$ws.Cells.Item(6, 10).Value = 0.2323568509805327
$ws.Cells.Item(6, 13).Value = 2.000037
$ws.Cells.Item(6, 14).Value = 6.000111
$ws.Cells.Item(6, 15).Value = 0.044360605600967
$ws.Cells.Item(6, 16).Value = 0.044360605600967
$ws.Cells.Item(6, 17).Value = 2.91616728143
$ws.Cells.Item(6, 18).Value = 26.24550553287
$ws.Cells.Item(6, 19).Value = 0.01030749062503008
$ws.Cells.Item(6, 20).Value = 0.01030749062503008
$ws.Cells.Item(7, 8).Value = 5.708772
$ws.Cells.Item(7, 9).Value = 0.3032511962008422
$ws.Cells.Item(7, 10).Value = 0.3032511962008422
$ws.Cells.Item(7, 13).Value = 29.01761566666667
$ws.Cells.Item(7, 14).Value = 87.052847
$ws.Cells.Item(7, 15).Value = 0.6436075952942075
$ws.Cells.Item(7, 16).Value = 0.6436075952942075
$ws.Cells.Item(7, 17).Value = 55.218317274876
$ws.Cells.Item(7, 18).Value = 496.964855473884
$ws.Cells.Item(7, 19).Value = 0.195174773156916
$ws.Cells.Item(7, 20).Value = 0.195174773156916
$ws.Cells.Item(8, 8).Value = 5.708772
$ws.Cells.Item(8, 9).Value = 0.3032511962008422
$ws.Cells.Item(8, 10).Value = 0.3032511962008422
$ws.Cells.Item(8, 15).Value = 0.05796185537580412
$ws.Cells.Item(8, 16).Value = 0.05796185537580412
$ws.Cells.Item(8, 17).Value = 4.972837709472
$ws.Cells.Item(8, 18).Value = 44.755539385248
$ws.Cells.Item(8, 19).Value = 0.01757700197673282
$ws.Cells.Item(8, 20).Value = 0.01757700197673282
$ws.Cells.Item(9, 8).Value = 5.708772
$ws.Cells.Item(9, 9).Value = 0.3032511962008422
$ws.Cells.Item(9, 10).Value = 0.3032511962008422
$ws.Cells.Item(9, 13).Value = 1.123006
$ws.Cells.Item(9, 14).Value = 3.369018
$ws.Cells.Item(9, 15).Value = 0.02490815232594174
$ws.Cells.Item(9, 16).Value = 0.02490815232594174
$ws.Cells.Item(9, 17).Value = 2.136995069544
$ws.Cells.Item(9, 18).Value = 19.232955625896
$ws.Cells.Item(9, 19).Value = 0.007553426987994624
$ws.Cells.Item(9, 20).Value = 0.007553426987994623
$ws.Cells.Item(10, 8).Value = 5.708772
$ws.Cells.Item(10, 9).Value = 0.3032511962008422
$ws.Cells.Item(10, 10).Value = 0.3032511962008422
$ws.Cells.Item(10, 13).Value = 10.33196133333333
$ws.Cells.Item(10, 14).Value = 30.995884
$ws.Cells.Item(10, 15).Value = 0.2291617914030796
$ws.Cells.Item(10, 16).Value = 0.2291617914030796
$ws.Cells.Item(10, 17).Value = 19.660937188272
$ws.Cells.Item(10, 18).Value = 176.948434694448
$ws.Cells.Item(10, 19).Value = 0.06949358736651176
$ws.Cells.Item(10, 20).Value = 0.06949358736651176
$ws.Cells.Item(11, 8).Value = 5.708772
$ws.Cells.Item(11, 9).Value = 0.3032511962008422
$ws.Cells.Item(11, 10).Value = 0.3032511962008422
$ws.Cells.Item(11, 13).Value = 2.000037
$ws.Cells.Item(11, 14).Value = 6.000111
$ws.Cells.Item(11, 15).Value = 0.044360605600967
$ws.Cells.Item(11, 16).Value = 0.044360605600967
$ws.Cells.Item(11, 17).Value = 3.805918408188
$ws.Cells.Item(11, 18).Value = 34.253265673692
$ws.Cells.Item(11, 19).Value = 0.01345240671268703
$ws.Cells.Item(11, 20).Value = 0.01345240671268703
$ws.Cells.Item(12, 7).Value = 2.914094333333333
$ws.Cells.Item(12, 8).Value = 8.742283
$ws.Cells.Item(12, 9).Value = 0.4643919528186251
$ws.Cells.Item(12, 10).Value = 0.4643919528186251
$ws.Cells.Item(12, 13).Value = 29.01761566666667
$ws.Cells.Item(12, 14).Value = 87.052847
$ws.Cells.Item(12, 15).Value = 0.6436075952942075
$ws.Cells.Item(12, 16).Value = 0.6436075952942075
$ws.Cells.Item(12, 17).Value = 84.56006938107789
$ws.Cells.Item(12, 18).Value = 761.040624429701
$ws.Cells.Item(12, 19).Value = 0.2988861880275763
$ws.Cells.Item(12, 20).Value = 0.2988861880275763
$ws.Cells.Item(13, 7).Value = 2.914094333333333
$ws.Cells.Item(13, 8).Value = 8.742283
$ws.Cells.Item(13, 9).Value = 0.4643919528186251
$ws.Cells.Item(13, 10).Value = 0.4643919528186251
$ws.Cells.Item(13, 15).Value = 0.05796185537580412
$ws.Cells.Item(13, 16).Value = 0.05796185537580412
$ws.Cells.Item(13, 17).Value = 7.615290042985778
$ws.Cells.Item(13, 18).Value = 68.537610386872
$ws.Cells.Item(13, 19).Value = 0.0269170192069604
$ws.Cells.Item(13, 20).Value = 0.0269170192069604
$ws.Cells.Item(14, 7).Value = 2.914094333333333
$ws.Cells.Item(14, 8).Value = 8.742283
$ws.Cells.Item(14, 9).Value = 0.4643919528186251
$ws.Cells.Item(14, 10).Value = 0.4643919528186251
$ws.Cells.Item(14, 13).Value = 1.123006
$ws.Cells.Item(14, 14).Value = 3.369018
$ws.Cells.Item(14, 15).Value = 0.02490815232594174
$ws.Cells.Item(14, 16).Value = 0.02490815232594174
$ws.Cells.Item(14, 17).Value = 3.272545420899333
$ws.Cells.Item(14, 18).Value = 29.452908788094
$ws.Cells.Item(14, 19).Value = 0.01156714549974786
$ws.Cells.Item(14, 20).Value = 0.01156714549974786
$ws.Cells.Item(15, 7).Value = 2.914094333333333
$ws.Cells.Item(15, 8).Value = 8.742283
$ws.Cells.Item(15, 9).Value = 0.4643919528186251
$ws.Cells.Item(15, 10).Value = 0.4643919528186251
$ws.Cells.Item(15, 13).Value = 10.33196133333333
$ws.Cells.Item(15, 14).Value = 30.995884
$ws.Cells.Item(15, 15).Value = 0.2291617914030796
$ws.Cells.Item(15, 16).Value = 0.2291617914030796
$ws.Cells.Item(15, 17).Value = 30.10830997368577
$ws.Cells.Item(15, 18).Value = 270.974789763172
$ws.Cells.Item(15, 19).Value = 0.1064208918210905
$ws.Cells.Item(15, 20).Value = 0.1064208918210905
$ws.Cells.Item(16, 7).Value = 2.914094333333333
$ws.Cells.Item(16, 8).Value = 8.742283
$ws.Cells.Item(16, 9).Value = 0.4643919528186251
$ws.Cells.Item(16, 10).Value = 0.4643919528186251
$ws.Cells.Item(16, 13).Value = 2.000037
$ws.Cells.Item(16, 14).Value = 6.000111
$ws.Cells.Item(16, 15).Value = 0.044360605600967
$ws.Cells.Item(16, 16).Value = 0.044360605600967
$ws.Cells.Item(16, 17).Value = 5.828296488157001
$ws.Cells.Item(16, 18).Value = 52.45466839341301
$ws.Cells.Item(16, 19).Value = 0.0206007082632499
$ws.Cells.Item(16, 20).Value = 0.0206007082632499
